$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 5786
$ws.Range("E2").Value = 494
$ws.Range("F2").Value = 494
$ws.Range("G2").Value = 457
$ws.Range("H2").Value = 376
$ws.Range("I2").Value = 374
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 10819
$ws.Range("L2").Value = 5935
$ws.Range("M2").Value = 4883
$ws.Range("N2").Value = 4881
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 385
$ws.Range("Q2").Value = 847
$ws.Range("R2").Value = -283
$ws.Range("S2").Value = 1005
$ws.Range("T2").Value = 296
$ws.Range("U2").Value = 552
$ws.Range("V2").Value = 4487
$ws.Range("W2").Value = 8.529999999999999
$ws.Range("X2").Value = 6.49
$ws.Range("Y2").Value = 8.07
$ws.Range("Z2").Value = 3.73
$ws.Range("AA2").Value = 121.54
$ws.Range("AB2").Value = 699.33
$ws.Range("AC2").Value = 4858
$ws.Range("AD2").Value = 18.88
$ws.Range("AE2").Value = 63404
$ws.Range("AF2").Value = 1.45
$ws.Range("AG2").Value = 750
$ws.Range("AH2").Value = 0.82
$ws.Range("AI2").Value = 15.43
$ws.Range("AJ2").Value = 7701936

$ws.Range("D3").Value = 5679
$ws.Range("E3").Value = 539
$ws.Range("F3").Value = 554
$ws.Range("G3").Value = 625
$ws.Range("H3").Value = 486
$ws.Range("I3").Value = 486
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10635
$ws.Range("L3").Value = 5098
$ws.Range("M3").Value = 5538
$ws.Range("N3").Value = 5536
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 402
$ws.Range("Q3").Value = 424
$ws.Range("R3").Value = 210
$ws.Range("S3").Value = -704
$ws.Range("T3").Value = 283
$ws.Range("U3").Value = 141
$ws.Range("V3").Value = 3726
$ws.Range("W3").Value = 9.48
$ws.Range("X3").Value = 8.550000000000001
$ws.Range("Y3").Value = 9.32
$ws.Range("Z3").Value = 4.53
$ws.Range("AA3").Value = 92.05
$ws.Range("AB3").Value = 881.05
$ws.Range("AC3").Value = 6187
$ws.Range("AD3").Value = 23.76
$ws.Range("AE3").Value = 68932
$ws.Range("AF3").Value = 2.13
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 0.68
$ws.Range("AI3").Value = 16.54
$ws.Range("AJ3").Value = 8035805

$ws.Range("O4").ClearContents()
$ws.Range("D4").Value = 5605
$ws.Range("E4").Value = 148
$ws.Range("F4").Value = 148
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 106
$ws.Range("I4").Value = 106
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 10666
$ws.Range("L4").Value = 4774
$ws.Range("M4").Value = 5892
$ws.Range("N4").Value = 5892
$ws.Range("P4").Value = 422
$ws.Range("Q4").Value = 19
$ws.Range("R4").Value = 458
$ws.Range("S4").Value = -647
$ws.Range("T4").Value = 206
$ws.Range("U4").Value = -187
$ws.Range("V4").Value = 2927
$ws.Range("W4").Value = 2.65
$ws.Range("X4").Value = 1.9
$ws.Range("Y4").Value = 1.86
$ws.Range("Z4").Value = 1
$ws.Range("AA4").Value = 81.02
$ws.Range("AB4").Value = 917.23
$ws.Range("AC4").Value = 1267
$ws.Range("AD4").Value = 83.27
$ws.Range("AE4").Value = 69818
$ws.Range("AF4").Value = 1.51
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 0.47
$ws.Range("AI4").Value = 39.7
$ws.Range("AJ4").Value = 8443868

$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("D5").Value = 5550
$ws.Range("E5").Value = 239
$ws.Range("F5").Value = 239
$ws.Range("G5").Value = -83
$ws.Range("H5").Value = -54
$ws.Range("I5").Value = -54
$ws.Range("K5").Value = 8997
$ws.Range("L5").Value = 3185
$ws.Range("M5").Value = 5812
$ws.Range("N5").Value = 5812
$ws.Range("P5").Value = 422
$ws.Range("Q5").Value = 512
$ws.Range("R5").Value = -21
$ws.Range("S5").Value = -1253
$ws.Range("T5").Value = 107
$ws.Range("U5").Value = 405
$ws.Range("V5").Value = 1803
$ws.Range("W5").Value = 4.3
$ws.Range("X5").Value = -0.97
$ws.Range("Y5").Value = -0.92
$ws.Range("Z5").Value = -0.55
$ws.Range("AA5").Value = 54.81
$ws.Range("AB5").Value = 900.72
$ws.Range("AC5").Value = -639
$ws.Range("AD5").Value = -155.87
$ws.Range("AE5").Value = 68875
$ws.Range("AF5").Value = 1.45
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 0.5
$ws.Range("AI5").Value = -78.2
$ws.Range("AJ5").Value = 8443868

$ws.Range("D6").Value = 5674
$ws.Range("E6").Value = 394
$ws.Range("F6").Value = 394
$ws.Range("G6").Value = 161
$ws.Range("H6").Value = 80
$ws.Range("I6").Value = 80
$ws.Range("K6").Value = 9487
$ws.Range("L6").Value = 3693
$ws.Range("M6").Value = 5794
$ws.Range("N6").Value = 5794
$ws.Range("P6").Value = 422
$ws.Range("Q6").Value = 501
$ws.Range("R6").Value = -442
$ws.Range("S6").Value = 192
$ws.Range("T6").Value = 77
$ws.Range("U6").Value = 424
$ws.Range("V6").Value = 2102
$ws.Range("W6").Value = 6.95
$ws.Range("X6").Value = 1.41
$ws.Range("Y6").Value = 1.38
$ws.Range("Z6").Value = 0.87
$ws.Range("AA6").Value = 63.74
$ws.Range("AB6").Value = 921.4400000000001
$ws.Range("AC6").Value = 948
$ws.Range("AD6").Value = 110.18
$ws.Range("AE6").Value = 68658
$ws.Range("AF6").Value = 1.52
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 0.96
$ws.Range("AI6").Value = 105.36
$ws.Range("AJ6").Value = 8443868

$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("AG7").Value = 1000
$ws.Range("AH7").Value = 1.06

$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AI8").ClearContents()
$ws.Range("AG8").Value = 1000
$ws.Range("AH8").Value = 1.06

$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AI9").ClearContents()
$ws.Range("AG9").Value = 1040
$ws.Range("AH9").Value = 1.11
